# Se agregan los script 0131/0132/0133/0134/0135 a la Clase Tests_AdmInstituciones
#
# The sheet currently ends with a short run of "template" rows (49-51, only
# styled/blank placeholders) followed by a couple of footer rows (52-54).
# We need to turn rows 49-51 into 15 new data rows (DEC_0132..DEC_0146,
# following exactly the same pattern used for DEC_0124..DEC_0131 in rows
# 41-48) while keeping the two footer rows intact further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three blank placeholder rows (old 49:51) - this shifts the old
# footer rows (52,53,54) up to become (49,50,51).
$ws.Range("A49:A51").EntireRow.Delete()

# Make room for the 15 new data rows by inserting 15 blank rows right
# before the (now shifted) footer rows - this pushes them back down to
# (64,65,66), matching the target layout.
$ws.Range("A49:A63").EntireRow.Insert()

# Fill the 15 new rows with the same shape as the preceding DEC_01xx rows:
#   A -> DEC_0132 .. DEC_0146
#   B -> "13712759-8"
#   C -> "Verity1.1"
#   D..J -> "SIN_DATO"
$codes = 132..146
$row = 49
foreach ($code in $codes) {
    $label = "DEC_0" + $code
    $ws.Cells.Item($row, 1).Value = $label
    $ws.Cells.Item($row, 2).Value = "13712759-8"
    $ws.Cells.Item($row, 3).Value = "Verity1.1"
    $ws.Cells.Item($row, 4).Value = "SIN_DATO"
    $ws.Cells.Item($row, 5).Value = "SIN_DATO"
    $ws.Cells.Item($row, 6).Value = "SIN_DATO"
    $ws.Cells.Item($row, 7).Value = "SIN_DATO"
    $ws.Cells.Item($row, 8).Value = "SIN_DATO"
    $ws.Cells.Item($row, 9).Value = "SIN_DATO"
    $ws.Cells.Item($row, 10).Value = "SIN_DATO"
    $row++
}

# Match the author's final cursor/viewport position.
$ws.Range("C62").Select()
